$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update subcategory (column H) text for rows whose label changed
$ws.Range("H2").Value = "photo(s)"
$ws.Range("H3").Value = "photo(s)"
$ws.Range("H4").Value = "photo(s)"
$ws.Range("H5").Value = "mixed statistical plot (more than 1 statistical plot and type)"
$ws.Range("H6").Value = "photo(s)"
$ws.Range("H7").Value = "mixed statistical plot (more than 1 statistical plot and type)"
$ws.Range("H10").Value = "photo(s)"
$ws.Range("H11").Value = "photo(s)"
$ws.Range("H13").Value = "photo(s)"
$ws.Range("H14").Value = "photo(s)"
$ws.Range("H15").Value = "photo(s)"
$ws.Range("H19").Value = "data display"
$ws.Range("H23").Value = "data collection, data analysis, data gathering diagram"
$ws.Range("H25").Value = "bar chart(s)"
$ws.Range("H30").Value = "bar chart(s)"
$ws.Range("H35").Value = "line graph(s)"
$ws.Range("H41").Value = "bar chart(s)"
$ws.Range("H45").Value = "line graph(s)"
$ws.Range("H46").Value = "line graph(s)"
$ws.Range("H47").Value = "line graph(s)"
$ws.Range("H48").Value = "line graph(s)"
$ws.Range("H50").Value = "photo(s)"
$ws.Range("H51").Value = "photo(s)"
$ws.Range("H56").Value = "data display"
$ws.Range("H66").Value = "photo(s)"
$ws.Range("H67").Value = "photo(s)"
$ws.Range("H68").Value = "data collection, data analysis, data gathering diagram"
$ws.Range("H74").Value = "line graph(s)"
$ws.Range("H77").Value = "line graph(s)"
$ws.Range("H94").Value = "data display"
$ws.Range("H104").Value = "line graph(s)"
$ws.Range("H105").Value = "line graph(s)"
$ws.Range("H107").Value = "line graph(s)"
$ws.Range("H108").Value = "line graph(s)"
$ws.Range("H113").Value = "photo(s)"
$ws.Range("H116").Value = "scatter plot(s)"
$ws.Range("H118").Value = "bar chart(s)"
$ws.Range("H121").Value = "line graph(s)"
$ws.Range("H128").Value = "data collection, data analysis, data gathering diagram"
$ws.Range("H135").Value = "line graph(s)"
$ws.Range("H137").Value = "line graph(s)"
$ws.Range("H138").Value = "line graph(s)"

# Remove the is_viewed column (I) entirely, shifting dimension to A1:H138
$ws.Columns.Item(9).Delete()

